$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column cells keep their exact textual representation
# (avoids Excel auto-converting numeric-looking strings to numbers,
# which would strip meaningful trailing zeros, e.g. "156.50" -> 156.5).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '63.824.05'
$ws.Range('E2').Value = '  +1.80%  '
$ws.Range('D3').Value = '3.409.47'
$ws.Range('E3').Value = '  +1.78%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '569.29'
$ws.Range('E5').Value = '  -0.82%  '
$ws.Range('D6').Value = '156.50'
$ws.Range('E6').Value = '  +2.67%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.410.00'
$ws.Range('E8').Value = '  +1.79%  '
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').Value = '  +2.83%  '
$ws.Range('D10').Value = '7.36'
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('E11').Value = '  +3.48%  '
$ws.Range('D12').Value = '0.433'
$ws.Range('E12').Value = '  -1.74%  '
$ws.Range('D13').Value = '3.999.86'
$ws.Range('E13').Value = '  +1.96%  '
$ws.Range('D14').Value = '0.134'
$ws.Range('E14').Value = '  -3.47%  '
$ws.Range('E15').Value = '  +6.53%  '
$ws.Range('D16').Value = '27.19'
$ws.Range('E16').Value = '  +1.12%  '
$ws.Range('D17').Value = '63.807.15'
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('D18').Value = '3.389.10'
$ws.Range('E18').Value = '  +2.53%  '
$ws.Range('E19').Value = '  -1.50%  '
$ws.Range('D20').Value = '14.04'
$ws.Range('E20').Value = '  +1.36%  '
$ws.Range('D21').Value = '381.34'
$ws.Range('E21').Value = '  -0.76%  '
$ws.Range('D22').Value = '8.08'
$ws.Range('E22').Value = '  -4.24%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '71.87'
$ws.Range('E24').Value = '  +2.57%  '
$ws.Range('D25').Value = '0.531'
$ws.Range('E25').Value = '  -0.83%  '
$ws.Range('E26').Value = '  +24.01%  '
$ws.Range('D27').Value = '9.36'
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  -0.40%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = '6.12'
$ws.Range('E30').Value = '  +8.44%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '1.36'
$ws.Range('E31').Value = '  +3.34%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '2.00'
$ws.Range('E32').Value = '  -0.74%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '6.43'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '23.20'
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '6.79'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('D37').Value = '160.47'
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').Value = '2.971.70'
$ws.Range('E39').Value = '  +6.99%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = '0.0757'
$ws.Range('E40').Value = '  +2.24%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '1.82'
$ws.Range('E41').Value = '  -3.18%  '
$ws.Range('D42').Value = '26.85'
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('D43').Value = '0.0313'
$ws.Range('E43').Value = '  -5.02%  '
$ws.Range('D44').Value = '41.86'
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('D45').Value = '0.758'
$ws.Range('E45').Value = '  +2.01%  '
$ws.Range('D46').Value = '4.29'
$ws.Range('E46').Value = '  +0.25%  '
$ws.Range('D47').Value = '23.18'
$ws.Range('E47').Value = '  +4.89%  '
$ws.Range('E48').Value = '  +2.94%  '
$ws.Range('D49').Value = '2.19'
$ws.Range('E49').Value = '  +20.96%  '
$ws.Range('D50').Value = '0.832'
$ws.Range('E50').Value = '  +2.88%  '
$ws.Range('D51').Value = '6.34'
$ws.Range('E51').Value = '  +0.10%  '
